$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 3035330.95
$ws.Range("C9").Value = 481606.99
$ws.Range("D9").Value = 3516937.94
$ws.Range("E9").Value = 13.69392915702118
$ws.Range("F9").Value = 86.30607084297883
$ws.Range("G9").Value = -53.45508581966081
$ws.Range("H9").Value = -45.1861185976823
$ws.Range("I9").Value = 30379
$ws.Range("J9").Value = 1296
$ws.Range("K9").Value = 31675
$ws.Range("L9").Value = 21831
$ws.Range("M9").Value = 161.0983436397783
$ws.Range("N9").Value = 9.985120410348269
